# semana 22 de 2024
# Insert a new surveillance event (298 - Evento adverso grave posterior a la
# vacunacion) ahead of event 300 and refresh the week's Esperado/Observado/
# valor-p figures for the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..36 down one position, opening up a blank row 10 for the
# new event, then fill it in.
$ws.Rows.Item(10).Insert()

# Force the event code into column A to stay text (matching every other
# "evento" code in the sheet) instead of Excel's default numeric coercion,
# then drop the temporary text format so no stray style is left behind.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "298"
$ws.Range("A10").Style = "Normal"

$ws.Range("B10").Value = "Evento adverso grave posterior a la vacunacion"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Refresh Esperado (C) / Observado (D) / valor p (E) figures for the rows
# that already existed, now at their post-insert row numbers.
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.14

$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 0.37

$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0.13

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 82

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.27

$ws.Range("C11").Value = 41
$ws.Range("D11").Value = 20

$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.37

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.15

$ws.Range("C15").Value = 12

$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 1

$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0.03

$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.14

$ws.Range("D24").Value = 1

$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0

$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0

$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 0.02

$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0.37

$ws.Range("D31").Value = 1

$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 1

$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0.37

$ws.Range("C34").Value = 7
$ws.Range("E34").Value = 0.05

$ws.Range("C35").Value = 10

$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 3
$ws.Range("E36").Value = 0.03
